$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two observations were dropped from the cleaned dataset: "RM 232" and "SC 92".
# Locate them by ID (column A) rather than hard-coded row numbers, then remove
# their entire rows (remaining rows shift up).
$rowRM232 = $ws.Columns.Item(1).Find("RM 232").Row
$ws.Rows.Item($rowRM232).Delete()

$rowSC92 = $ws.Columns.Item(1).Find("SC 92").Row
$ws.Rows.Item($rowSC92).Delete()

# The imputed/missing values in column C ("D" header, col D) were re-rolled:
# some previously-missing cells now have a value, and some previously-filled
# cells are now missing.
$ws.Range("D3").Value = -14.2     # RM 8: was missing, now has a value
$ws.Range("D5").Value = ""        # RM 14: was -14.4, now missing
$ws.Range("D21").Value = -14.3    # RM 135: was missing, now has a value
$ws.Range("D23").Value = ""       # RM 140: was -13.9, now missing
$ws.Range("D32").Value = -14.7    # SC 193: was missing, now has a value
